$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 479
$ws.Range("F7").Value = 539
$ws.Range("F11").Value = 766
$ws.Range("F12").Value = 381
$ws.Range("F14").Value = 386
$ws.Range("F16").Value = 1019
$ws.Range("F17").Value = 19715
$ws.Range("F18").Value = 582
$ws.Range("F19").Value = 60
$ws.Range("F20").Value = 219
$ws.Range("F24").Value = 13
$ws.Range("F28").Value = 325

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 201
$ws.Range("F7").Value = 221
$ws.Range("F8").Value = 3384
$ws.Range("F10").Value = 87
$ws.Range("F16").Value = 3023

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 89
$ws.Range("F4").Value = 559

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 89
$ws.Range("F6").Value = 559
$ws.Range("F7").Value = 479
$ws.Range("F11").Value = 539
$ws.Range("F16").Value = 201
$ws.Range("F20").Value = 766
$ws.Range("F21").Value = 381
$ws.Range("F23").Value = 386
$ws.Range("F25").Value = 1019
$ws.Range("F26").Value = 19715
$ws.Range("F27").Value = 221
$ws.Range("F28").Value = 3384
$ws.Range("F30").Value = 87
$ws.Range("F32").Value = 582
$ws.Range("F33").Value = 61
$ws.Range("F34").Value = 219
$ws.Range("F40").Value = 13
$ws.Range("F46").Value = 325
$ws.Range("F48").Value = 3024
